$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.960.49"
$ws.Range("E2").Value = "  -1.55%  "

# Row 3
$ws.Range("D3").Value = "3.949.68"
$ws.Range("E3").Value = "  -2.06%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'536.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.98%  "

# Row 6
$ws.Range("D6").Value = "'148.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "

# Row 7
$ws.Range("D7").Value = "3.944.88"

# Row 8
$ws.Range("E8").Value = "  -5.61%  "

# Row 9
$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").Value = "'0.736"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.91%  "

# Row 11
$ws.Range("E11").Value = "  -5.44%  "

# Row 12
$ws.Range("D12").Value = "'55.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +14.64%  "

# Row 13
$ws.Range("D13").Value = "'0.0000315"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.69%  "

# Row 14
$ws.Range("D14").Value = "'10.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.97%  "

# Row 15
$ws.Range("D15").Value = "4.587.26"
$ws.Range("E15").Value = "  -1.97%  "

# Row 16
$ws.Range("D16").Value = "3.963.92"
$ws.Range("E16").Value = "  -2.36%  "

# Row 17
$ws.Range("D17").Value = "'20.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.34%  "

# Row 19
$ws.Range("E19").Value = "  -1.54%  "

# Row 20
$ws.Range("D20").Value = "'1.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.25%  "

# Row 21
$ws.Range("D21").Value = "70.857.08"
$ws.Range("E21").Value = "  -1.68%  "

# Row 22
$ws.Range("D22").Value = "'424.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.31%  "

# Row 23
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("D24").Value = "'96.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.59%  "

# Row 25
$ws.Range("D25").Value = "'4.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.45%  "

# Row 26
$ws.Range("D26").Value = "'14.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.61%  "

# Row 27
$ws.Range("D27").Value = "'11.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.20%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'3.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.68%  "

# Row 29
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'10.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.13%  "

# Row 30
$ws.Range("D30").Value = "'5.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.31%  "

# Row 31
$ws.Range("D31").Value = "'36.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.05%  "

# Row 32
$ws.Range("D32").Value = "'7.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +15.62%  "

# Row 33
$ws.Range("D33").Value = "'50.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +19.76%  "

# Row 34
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'13.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.99%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.15%  "

# Row 36
$ws.Range("D36").Value = "'683.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.04%  "

# Row 37
$ws.Range("D37").Value = "'64.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.24%  "

# Row 38
$ws.Range("D38").Value = "'0.436"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.34%  "

# Row 39
$ws.Range("E39").Value = "  -5.11%  "

# Row 40
$ws.Range("D40").Value = "'0.149"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.92%  "

# Row 41
$ws.Range("D41").Value = "'3.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.76%  "

# Row 42
$ws.Range("E42").Value = "  +0.05%  "

# Row 43
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.18%  "

# Row 44
$ws.Range("D44").Value = "'0.0479"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.42%  "

# Row 45
$ws.Range("E45").Value = "  -2.31%  "

# Row 46
$ws.Range("E46").Value = "  -7.17%  "

# Row 47
$ws.Range("D47").Value = "'9.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.89%  "

# Row 48
$ws.Range("D48").Value = "'2.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.81%  "

# Row 49
$ws.Range("D49").Value = "'3.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.15%  "

# Row 50
$ws.Range("E50").Value = "  -2.12%  "

# Row 51
$ws.Range("D51").Value = "'0.000271"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.64%  "
